# AoC 2024 - Day 24 (solution, part I)
# Adds a new "Day 24" worksheet after the last existing sheet ("Day 20")
# and fills in the half-adder / full-adder truth-table sketch.

$wb = $excel.ActiveWorkbook

# --- Add the new sheet at the very end of the tab strip -------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Day 24"

# --- Bit-weight header row (bold) ------------------------------------------
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 0
$ws.Range("B1:D1").Font.Bold = $true

# --- x / y binary rows -------------------------------------------------
$ws.Range("A2").Value = "x"
$ws.Range("A2").Font.Bold = $true
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("G2").Value = 2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1

$ws.Range("A3").Value = "y"
$ws.Range("A3").Font.Bold = $true
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 3

# --- small scratch area -----------------------------------------------
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0

# --- half adder truth table (filled row by row; this also establishes
#     the shared-string insertion order: "-->" before "s" before "c",
#     and "XOR" before "AND") ------------------------------------------
$ws.Range("A12").Value = 0
$ws.Range("B12").Value = "+"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = "'-->"
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0

$ws.Range("A13").Value = 0
$ws.Range("B13").Value = "+"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "'-->"
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "+"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = "'-->"
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1

$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "+"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "'-->"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0

# --- half adder truth table header (filled after the body, matching
#     the original author's shared-string order) ------------------------
$ws.Range("F11").Value = "s"
$ws.Range("E11").Value = "c"

# --- gate legend ---------------------------------------------------------
$ws.Range("F16").Value = "XOR"
$ws.Range("E16").Value = "AND"

# --- page setup: portrait orientation ------------------------------------
$ws.PageSetup.Orientation = 1

# --- view state: zoom, selection, active tab ------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 175
$ws.Range("K3").Select()
